$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure target cells keep their literal text representation (avoid Excel
# auto-converting numeric-looking strings like "308.21" or "0.85%" into
# real numbers / percentages).
$cells = @(
    "D2"
    "E2"
    "D3"
    "E3"
    "B4"
    "C4"
    "D4"
    "E4"
    "B5"
    "C5"
    "D5"
    "E5"
    "B6"
    "C6"
    "D6"
    "E6"
    "B7"
    "C7"
    "D7"
    "E7"
    "B8"
    "C8"
    "D8"
    "E8"
    "B9"
    "C9"
    "D9"
    "E9"
    "B10"
    "C10"
    "D10"
    "E10"
    "B11"
    "C11"
    "D11"
    "E11"
    "B12"
    "C12"
    "D12"
    "E12"
    "B13"
    "C13"
    "D13"
    "E13"
    "B14"
    "C14"
    "D14"
    "E14"
    "B15"
    "C15"
    "D15"
    "E15"
    "B16"
    "C16"
    "D16"
    "E16"
    "B17"
    "C17"
    "D17"
    "E17"
    "E18"
    "D19"
    "E19"
    "D20"
    "E20"
    "D21"
    "E21"
    "D22"
    "E22"
    "D23"
    "E23"
    "D24"
    "E24"
    "D25"
    "E25"
    "E26"
    "D27"
    "E27"
    "D39"
    "E39"
    "D40"
    "E40"
    "D41"
    "E41"
    "D42"
    "E42"
    "D43"
    "E43"
    "D44"
    "E44"
    "D45"
    "E45"
    "E46"
    "E47"
    "D48"
    "E49"
    "D50"
    "E50"
    "E51"
)
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '308.21'
$ws.Range("E2").Value = '0.85%'
$ws.Range("D3").Value = '38.84'
$ws.Range("E3").Value = '8.77%'
$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").Value = '5.105'
$ws.Range("E4").Value = '1.23%'
$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").Value = '0.08141'
$ws.Range("E5").Value = '1.38%'
$ws.Range("B6").Value = 'FTXToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D6").Value = '1.961'
$ws.Range("E6").Value = '2.42%'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").Value = '7.929'
$ws.Range("E7").Value = '1.96%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").Value = '0.9273'
$ws.Range("E8").Value = '0.52%'
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").Value = '0.1442'
$ws.Range("E9").Value = '12.26%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1951'
$ws.Range("E10").Value = '1.73%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = '0.09258'
$ws.Range("E11").Value = '1.27%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").Value = '0.03509'
$ws.Range("E12").Value = '2.03%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").Value = '0.09845'
$ws.Range("E13").Value = '-0.05%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").Value = '0.001402'
$ws.Range("E14").Value = '0.09%'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006034'
$ws.Range("E15").Value = '-4.13%'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").Value = '3.601'
$ws.Range("E16").Value = '-3.25%'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").Value = '4.202'
$ws.Range("E17").Value = '1.13%'
$ws.Range("E18").Value = '2.37%'
$ws.Range("D19").Value = '0.3446'
$ws.Range("E19").Value = '0.04%'
$ws.Range("D20").Value = '0.1334'
$ws.Range("E20").Value = '-0.51%'
$ws.Range("D21").Value = '4.821'
$ws.Range("E21").Value = '-6.70%'
$ws.Range("D22").Value = '0.2411'
$ws.Range("E22").Value = '-7.30%'
$ws.Range("D23").Value = '0.04447'
$ws.Range("E23").Value = '0.46%'
$ws.Range("D24").Value = '0.001244'
$ws.Range("E24").Value = '0.75%'
$ws.Range("D25").Value = '0.004854'
$ws.Range("E25").Value = '4.94%'
$ws.Range("E26").Value = '0.00%'
$ws.Range("D27").Value = '0.0001303'
$ws.Range("E27").Value = '4.08%'
$ws.Range("D39").Value = '0.02102'
$ws.Range("E39").Value = '7.39%'
$ws.Range("D40").Value = '0.05110'
$ws.Range("E40").Value = '-8.19%'
$ws.Range("D41").Value = '0.007472'
$ws.Range("E41").Value = '-2.15%'
$ws.Range("D42").Value = '0.01013'
$ws.Range("E42").Value = '-0.58%'
$ws.Range("D43").Value = '0.1365'
$ws.Range("E43").Value = '0.84%'
$ws.Range("D44").Value = '0.002145'
$ws.Range("E44").Value = '-0.85%'
$ws.Range("D45").Value = '0.01051'
$ws.Range("E45").Value = '5.85%'
$ws.Range("E46").Value = '1.49%'
$ws.Range("E47").Value = '0.15%'
$ws.Range("D48").Value = '0.003070'
$ws.Range("E49").Value = '-3.36%'
$ws.Range("D50").Value = '0.00002104'
$ws.Range("E50").Value = '0.15%'
$ws.Range("E51").Value = '0.15%'
